$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 159.36363
$ws.Range("J55").Value = 223.5
$ws.Range("L55").Value = 223.5
$ws.Range("N55").Value = -651.5
$ws.Range("H125").Value = 10192394
$ws.Range("I125").Value = 599
$ws.Range("J125").Value = 11211573
$ws.Range("K125").Value = 5391
$ws.Range("L125").Value = 100904157
$ws.Range("M125").Value = -2931
$ws.Range("N125").Value = -100909077
$ws.Range("H131").Value = 5228.269
$ws.Range("I131").Value = 1101.9445
$ws.Range("J131").Value = 14512.5
$ws.Range("K131").Value = 3305.8335
$ws.Range("L131").Value = 43537.5
$ws.Range("M131").Value = 1734.1665
$ws.Range("N131").Value = -53617.5
$ws.Range("H132").Value = 278737.28
$ws.Range("I132").Value = 329166.8
$ws.Range("J132").Value = 45500.75
$ws.Range("K132").Value = 987500.3999999999
$ws.Range("L132").Value = 136502.25
$ws.Range("M132").Value = -984970.3999999999
$ws.Range("N132").Value = -141562.25
$ws.Range("H137").Value = 1730.0952
$ws.Range("I137").Value = 1075.3334
$ws.Range("J137").Value = 1992
$ws.Range("K137").Value = 3226.0002
$ws.Range("L137").Value = 5976
$ws.Range("M137").Value = -676.0001999999999
$ws.Range("N137").Value = -11076
$ws.Range("H138").Value = 5714632
$ws.Range("I138").Value = 920880.0600000001
$ws.Range("J138").Value = 8477133
$ws.Range("K138").Value = 2762640.18
$ws.Range("L138").Value = 25431399
$ws.Range("M138").Value = -2757500.18
$ws.Range("N138").Value = -25441679

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2192.037
$ws.Range("I61").Value = 1136.45
$ws.Range("K61").Value = 1136.45
$ws.Range("M61").Value = -924.45
$ws.Range("H136").Value = 2192.037
$ws.Range("I136").Value = 1136.45
$ws.Range("K136").Value = 3409.35
$ws.Range("M136").Value = -859.3500000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3285.9656
$ws.Range("I134").Value = 2274.4092
$ws.Range("J134").Value = 6465.143
$ws.Range("K134").Value = 6823.2276
$ws.Range("L134").Value = 19395.429
$ws.Range("M134").Value = -4288.2276
$ws.Range("N134").Value = -24465.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2126.8684
$ws.Range("I31").Value = 1276.28
$ws.Range("J31").Value = 3762.6155
$ws.Range("K31").Value = 1276.28
$ws.Range("L31").Value = 3762.6155
$ws.Range("M31").Value = -981.28
$ws.Range("N31").Value = -4352.6155
$ws.Range("H34").Value = 2126.8684
$ws.Range("I34").Value = 1276.28
$ws.Range("J34").Value = 3762.6155
$ws.Range("K34").Value = 1276.28
$ws.Range("L34").Value = 3762.6155
$ws.Range("M34").Value = -1074.28
$ws.Range("N34").Value = -4166.6155
$ws.Range("H122").Value = 1390.2222
$ws.Range("I122").Value = 1022.4
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 3067.2
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -617.1999999999998
$ws.Range("N122").Value = -10450
$ws.Range("H132").Value = 2869.9583
$ws.Range("I132").Value = 1803.125
$ws.Range("J132").Value = 5003.625
$ws.Range("K132").Value = 5409.375
$ws.Range("L132").Value = 15010.875
$ws.Range("M132").Value = -2879.375
$ws.Range("N132").Value = -20070.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1604.8
$ws.Range("I64").Value = 674.6667
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 2024.0001
$ws.Range("L64").Value = 9000
$ws.Range("M64").Value = -1754.0001
$ws.Range("N64").Value = -9540
$ws.Range("H67").Value = 1604.8
$ws.Range("I67").Value = 674.6667
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 2024.0001
$ws.Range("L67").Value = 9000
$ws.Range("M67").Value = -1088.0001
$ws.Range("N67").Value = -10872
$ws.Range("H76").Value = 3666.6667
$ws.Range("I76").Value = 1000
$ws.Range("J76").Value = 9000
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 27000
$ws.Range("M76").Value = -2617
$ws.Range("N76").Value = -27766
$ws.Range("H79").Value = 3666.6667
$ws.Range("I79").Value = 1000
$ws.Range("J79").Value = 9000
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 27000
$ws.Range("M79").Value = -1674
$ws.Range("N79").Value = -29652
$ws.Range("H94").Value = 2500
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 7500
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -8852
$ws.Range("H120").Value = 67515
$ws.Range("I120").Value = 100030
$ws.Range("K120").Value = 300090
$ws.Range("M120").Value = -295252
$ws.Range("H125").Value = 2818.5
$ws.Range("J125").Value = 2995.2144
$ws.Range("L125").Value = 8985.643199999999
$ws.Range("N125").Value = -18825.6432
$ws.Range("H131").Value = 2964.4688
$ws.Range("J131").Value = 3100.4333
$ws.Range("L131").Value = 9301.2999
$ws.Range("N131").Value = -19381.2999
$ws.Range("H136").Value = 12090.5
$ws.Range("I136").Value = 2010
$ws.Range("J136").Value = 16410.715
$ws.Range("K136").Value = 6030
$ws.Range("L136").Value = 49232.145
$ws.Range("M136").Value = -930
$ws.Range("N136").Value = -59432.145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1944.0883
$ws.Range("I126").Value = 1387.4375
$ws.Range("J126").Value = 2438.889
$ws.Range("K126").Value = 4162.3125
$ws.Range("L126").Value = 7316.667
$ws.Range("M126").Value = -1692.3125
$ws.Range("N126").Value = -12256.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 50000
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H68").Value = 2114.5715
$ws.Range("I68").Value = 1800.3334
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 1800.3334
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1051.3334
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 2114.5715
$ws.Range("I71").Value = 1800.3334
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 9001.666999999999
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -5257.666999999999
$ws.Range("N71").Value = -27488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 36666.668
$ws.Range("J28").Value = 30000
$ws.Range("L28").Value = 30000
$ws.Range("N28").Value = -30696
$ws.Range("H107").Value = 2778398.5
$ws.Range("I107").Value = 3268603.5
$ws.Range("J107").Value = 570.3333
$ws.Range("K107").Value = 9805810.5
$ws.Range("L107").Value = 1710.9999
$ws.Range("M107").Value = -9803890.5
$ws.Range("N107").Value = -5550.9999
$ws.Range("H122").Value = 63687.688
$ws.Range("I122").Value = 112089.78
$ws.Range("J122").Value = 1456.4286
$ws.Range("K122").Value = 336269.34
$ws.Range("L122").Value = 4369.2858
$ws.Range("M122").Value = -333819.34
$ws.Range("N122").Value = -9269.2858
$ws.Range("H132").Value = 10639878
$ws.Range("I132").Value = 14707094
$ws.Range("K132").Value = 44121282
$ws.Range("M132").Value = -44118752
